# Swap the data contents of row 4 and row 5 for the columns that
# actually differ between the two records (A, B, E, F, G, H, I, J, Q, R).
# All other columns (C, D, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG,
# AT, AW, AX, AY, ...) already hold identical values in both rows, so
# they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "E", "F", "G", "H", "I", "J", "Q", "R")

foreach ($col in $columns) {
    $addr4 = "$col" + "4"
    $addr5 = "$col" + "5"

    $val4 = $ws.Range($addr4).Value2
    $val5 = $ws.Range($addr5).Value2

    $ws.Range($addr4).Value2 = $val5
    $ws.Range($addr5).Value2 = $val4
}
